$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 143.33333

$ws.Range("H33").Value = 362.57895
$ws.Range("I33").Value = 319.93332
$ws.Range("J33").Value = 522.5
$ws.Range("K33").Value = 319.93332
$ws.Range("L33").Value = 522.5
$ws.Range("M33").Value = -90.93331999999998
$ws.Range("N33").Value = -980.5

$ws.Range("H64").Value = 3501.111
$ws.Range("I64").Value = 3297.9167
$ws.Range("J64").Value = 3663.6667
$ws.Range("K64").Value = 3297.9167
$ws.Range("L64").Value = 3663.6667
$ws.Range("M64").Value = -3049.9167
$ws.Range("N64").Value = -4159.6667

$ws.Range("H67").Value = 3501.111
$ws.Range("I67").Value = 3297.9167
$ws.Range("J67").Value = 3663.6667
$ws.Range("K67").Value = 3297.9167
$ws.Range("L67").Value = 3663.6667
$ws.Range("M67").Value = -2439.9167
$ws.Range("N67").Value = -5379.6667

$ws.Range("H113").Value = 2956.111
$ws.Range("I113").Value = 2401.6667
$ws.Range("J113").Value = 3233.3333
$ws.Range("K113").Value = 2401.6667
$ws.Range("L113").Value = 3233.3333
$ws.Range("M113").Value = 852.3332999999998
$ws.Range("N113").Value = -9741.3333

$ws.Range("H116").Value = 2252.6086
$ws.Range("I116").Value = 1410
$ws.Range("J116").Value = 2550
$ws.Range("K116").Value = 1410
$ws.Range("L116").Value = 2550
$ws.Range("M116").Value = 2032
$ws.Range("N116").Value = -9434

$ws.Range("H132").Value = 2913076.5
$ws.Range("I132").Value = 3478763.2
$ws.Range("J132").Value = 3830.2856
$ws.Range("K132").Value = 10436289.6
$ws.Range("L132").Value = 11490.8568
$ws.Range("M132").Value = -10433759.6
$ws.Range("N132").Value = -16550.8568

$ws.Range("H137").Value = 2760.32
$ws.Range("I137").Value = 2715.0293
$ws.Range("J137").Value = 2856.5625
$ws.Range("K137").Value = 8145.0879
$ws.Range("L137").Value = 8569.6875
$ws.Range("M137").Value = -5595.0879
$ws.Range("N137").Value = -13669.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1031
$ws.Range("I2").Value = 1051.875
$ws.Range("J2").Value = 919.6667
$ws.Range("K2").Value = 1051.875
$ws.Range("L2").Value = 919.6667
$ws.Range("M2").Value = -938.875
$ws.Range("N2").Value = -1145.6667

$ws.Range("H74").Value = 5343.1924
$ws.Range("I74").Value = 5371.857
$ws.Range("J74").Value = 5222.8
$ws.Range("K74").Value = 5371.857
$ws.Range("L74").Value = 5222.8
$ws.Range("M74").Value = -4497.857
$ws.Range("N74").Value = -6970.8

$ws.Range("H77").Value = 5343.1924
$ws.Range("I77").Value = 5371.857
$ws.Range("J77").Value = 5222.8
$ws.Range("K77").Value = 26859.285
$ws.Range("L77").Value = 26114
$ws.Range("M77").Value = -22491.285
$ws.Range("N77").Value = -34850

$ws.Range("H102").Value = 5735.9
$ws.Range("I102").Value = 4959.8
$ws.Range("J102").Value = 6512
$ws.Range("K102").Value = 4959.8
$ws.Range("L102").Value = 6512
$ws.Range("M102").Value = -3337.8
$ws.Range("N102").Value = -9756

$ws.Range("H116").Value = 1031
$ws.Range("I116").Value = 1051.875
$ws.Range("J116").Value = 919.6667
$ws.Range("K116").Value = 1051.875
$ws.Range("L116").Value = 919.6667
$ws.Range("M116").Value = 1242.125
$ws.Range("N116").Value = -5507.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1031
$ws.Range("I3").Value = 1051.875
$ws.Range("J3").Value = 919.6667
$ws.Range("K3").Value = 1051.875
$ws.Range("L3").Value = 919.6667
$ws.Range("M3").Value = -937.875
$ws.Range("N3").Value = -1147.6667

$ws.Range("H86").Value = 2889.7576
$ws.Range("I86").Value = 3013.6296
$ws.Range("J86").Value = 2332.3333
$ws.Range("K86").Value = 3013.6296
$ws.Range("L86").Value = 2332.3333
$ws.Range("M86").Value = -1890.6296
$ws.Range("N86").Value = -4578.3333

$ws.Range("H89").Value = 2889.7576
$ws.Range("I89").Value = 3013.6296
$ws.Range("J89").Value = 2332.3333
$ws.Range("K89").Value = 15068.148
$ws.Range("L89").Value = 11661.6665
$ws.Range("M89").Value = -9452.148000000001
$ws.Range("N89").Value = -22893.6665

$ws.Range("H99").Value = 1751.25
$ws.Range("I99").Value = 1742
$ws.Range("J99").Value = 1766.6666
$ws.Range("K99").Value = 1742
$ws.Range("L99").Value = 1766.6666
$ws.Range("M99").Value = -244
$ws.Range("N99").Value = -4762.6666

$ws.Range("H132").Value = 53780
$ws.Range("J132").Value = 53780
$ws.Range("L132").Value = 53780
$ws.Range("N132").Value = -63900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25765.268
$ws.Range("I31").Value = 37325.965
$ws.Range("J31").Value = 2643.8667
$ws.Range("K31").Value = 37325.965
$ws.Range("L31").Value = 2643.8667
$ws.Range("M31").Value = -37030.965
$ws.Range("N31").Value = -3233.8667

$ws.Range("H34").Value = 25765.268
$ws.Range("I34").Value = 37325.965
$ws.Range("J34").Value = 2643.8667
$ws.Range("K34").Value = 37325.965
$ws.Range("L34").Value = 2643.8667
$ws.Range("M34").Value = -37123.965
$ws.Range("N34").Value = -3047.8667

$ws.Range("H86").Value = 3312.125
$ws.Range("I86").Value = 2874.25
$ws.Range("J86").Value = 3750
$ws.Range("K86").Value = 2874.25
$ws.Range("L86").Value = 3750
$ws.Range("M86").Value = -1751.25
$ws.Range("N86").Value = -5996

$ws.Range("H89").Value = 3312.125
$ws.Range("I89").Value = 2874.25
$ws.Range("J89").Value = 3750
$ws.Range("K89").Value = 14371.25
$ws.Range("L89").Value = 18750
$ws.Range("M89").Value = -8755.25
$ws.Range("N89").Value = -29982

$ws.Range("H134").Value = 13689.765
$ws.Range("I134").Value = 12368
$ws.Range("J134").Value = 15176.75
$ws.Range("K134").Value = 37104
$ws.Range("L134").Value = 45530.25
$ws.Range("M134").Value = -34569
$ws.Range("N134").Value = -50600.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 2630.8096
$ws.Range("I94").Value = 1391.25
$ws.Range("J94").Value = 2922.4707
$ws.Range("K94").Value = 4173.75
$ws.Range("L94").Value = 8767.4121
$ws.Range("M94").Value = -3497.75
$ws.Range("N94").Value = -10119.4121

$ws.Range("H131").Value = 939044.25
$ws.Range("I131").Value = 215
$ws.Range("J131").Value = 987189.4
$ws.Range("K131").Value = 645
$ws.Range("L131").Value = 2961568.2
$ws.Range("M131").Value = 4395
$ws.Range("N131").Value = -2971648.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2871.25
$ws.Range("I80").Value = 2619.375
$ws.Range("K80").Value = 2619.375
$ws.Range("M80").Value = -1621.375

$ws.Range("H83").Value = 2871.25
$ws.Range("I83").Value = 2619.375
$ws.Range("K83").Value = 13096.875
$ws.Range("M83").Value = -8104.875

$ws.Range("H97").Value = 753.7692
$ws.Range("I97").Value = 753.7692
$ws.Range("K97").Value = 753.7692
$ws.Range("M97").Value = -257.7692

$ws.Range("H113").Value = 3609.25
$ws.Range("I113").Value = 3623.4443
$ws.Range("J113").Value = 3566.6667
$ws.Range("K113").Value = 3623.4443
$ws.Range("L113").Value = 3566.6667
$ws.Range("M113").Value = -1453.4443
$ws.Range("N113").Value = -7906.6667

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3811.0476
$ws.Range("I122").Value = 4112.8887
$ws.Range("K122").Value = 12338.6661
$ws.Range("M122").Value = -9888.666100000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3354578
$ws.Range("I122").Value = 1814921.1
$ws.Range("J122").Value = 6947110.5
$ws.Range("K122").Value = 5444763.300000001
$ws.Range("L122").Value = 20841331.5
$ws.Range("M122").Value = -5442313.300000001
$ws.Range("N122").Value = -20846231.5
